$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen all 68 data columns (A:BP) from the narrow 2.16796875 char width to 12.7109375
for ($c = 1; $c -le 68; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 12.7109375
}

# Update connectivity-matrix cell values (restoring original weighted values)
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 0.5973745041742993
$ws.Range("E1").Value = 0.5431486480370014
$ws.Range("V1").Value = 0.6785952691944137
$ws.Range("AM1").Value = 0.827977111913237
$ws.Range("AS1").Value = 0
$ws.Range("BB1").Value = 0.5612631618123902
$ws.Range("BJ1").Value = 0.8634118042535076
$ws.Range("A2").Value = 0
$ws.Range("M2").Value = 0.6293636911876017
$ws.Range("N2").Value = 0
$ws.Range("AX2").Value = 0
$ws.Range("BB2").Value = 0
$ws.Range("BG2").Value = 0
$ws.Range("BO2").Value = 0.8716787027593691
$ws.Range("BP2").Value = 0
$ws.Range("A3").Value = 0.7431360849920408
$ws.Range("I3").Value = 0
$ws.Range("N3").Value = 0.5193282381099664
$ws.Range("AA3").Value = 0
$ws.Range("AZ3").Value = 0.844105669781331
$ws.Range("BB3").Value = 0.7139476303911625
$ws.Range("BK3").Value = 0
$ws.Range("E4").Value = 0.9196742796185104
$ws.Range("K4").Value = 0.8019854178228141
$ws.Range("L4").Value = 0
$ws.Range("V4").Value = 0.5603501675434281
$ws.Range("AB4").Value = 0
$ws.Range("AW4").Value = 0.7864265134246692
$ws.Range("A5").Value = 0.5653658094838034
$ws.Range("D5").Value = 0.9911516907489086
$ws.Range("J5").Value = 0
$ws.Range("U5").Value = 0.7660671578639239
$ws.Range("BK5").Value = 0.620968522722405
$ws.Range("BL5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("AN6").Value = 0
$ws.Range("BA6").Value = 0
$ws.Range("BF6").Value = 0.7905241873372741
$ws.Range("BH6").Value = 0.8980253870152032
$ws.Range("BO6").Value = 0.9950699479247356
$ws.Range("K7").Value = 0.7428160416726115
$ws.Range("N7").Value = 0
$ws.Range("V7").Value = 0
$ws.Range("AB7").Value = 0
$ws.Range("AI7").Value = 0
$ws.Range("AY7").Value = 0.6170983683527096
$ws.Range("F8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 0.5949224056804061
$ws.Range("P8").Value = 0
$ws.Range("V8").Value = 0.7626161511034777
$ws.Range("AO8").Value = 0
$ws.Range("AR8").Value = 0.5756028717847699
$ws.Range("BG8").Value = 0.5438819968521433
$ws.Range("C9").Value = 0
$ws.Range("M9").Value = 0.7017744676446672
$ws.Range("AM9").Value = 0
$ws.Range("AV9").Value = 0.5903157260785428
$ws.Range("AX9").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("L10").Value = 0.8147437788380062
$ws.Range("V10").Value = 0
$ws.Range("Y10").Value = 0.5374038223110442
$ws.Range("AA10").Value = 0.5150649407818851
$ws.Range("AU10").Value = 0.6342684390259835
$ws.Range("BI10").Value = 0
$ws.Range("BO10").Value = 0
$ws.Range("D11").Value = 0.6608828761548479
$ws.Range("G11").Value = 0.8694081018262158
$ws.Range("H11").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("X11").Value = 0
$ws.Range("AK11").Value = 0.9181129768392498
$ws.Range("AN11").Value = 0
$ws.Range("AT11").Value = 0
$ws.Range("AZ11").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("H12").Value = 0.8512364709508878
$ws.Range("J12").Value = 0.8330522142437401
$ws.Range("AH12").Value = 0
$ws.Range("AL12").Value = 0.6250814576594175
$ws.Range("AV12").Value = 0.5995604098170151
$ws.Range("AY12").Value = 0
$ws.Range("BD12").Value = 0
$ws.Range("BN12").Value = 0
$ws.Range("B13").Value = 0.8000832377733686
$ws.Range("I13").Value = 0.657354135965702
$ws.Range("O13").Value = 0
$ws.Range("Q13").Value = 0.7563259057308249
$ws.Range("AJ13").Value = 0.5298024096555617
$ws.Range("AK13").Value = 0.712165967548545
$ws.Range("AL13").Value = 0
$ws.Range("AT13").Value = 0.5586983527733913
$ws.Range("AU13").Value = 0
$ws.Range("AZ13").Value = 0
$ws.Range("BB13").Value = 0
$ws.Range("BP13").Value = 0.634437916570324
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0.812374505207652
$ws.Range("G14").Value = 0
$ws.Range("Q14").Value = 0.8398644510227629
$ws.Range("AA14").Value = 0.9409554370955941
$ws.Range("AD14").Value = 0.9051544707556759
$ws.Range("AN14").Value = 0.9060225411541579
$ws.Range("BF14").Value = 0.6186742655166879
$ws.Range("M15").Value = 0
$ws.Range("AD15").Value = 0
$ws.Range("AP15").Value = 0.9045146955778328
$ws.Range("AW15").Value = 0.7136451930783231
$ws.Range("BF15").Value = 0
$ws.Range("BJ15").Value = 0
$ws.Range("BN15").Value = 0
$ws.Range("BO15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("Q16").Value = 0.697335747413308
$ws.Range("AC16").Value = 0.889181224958684
$ws.Range("AD16").Value = 0.7313511871455787
$ws.Range("BG16").Value = 0
$ws.Range("BN16").Value = 0
$ws.Range("BP16").Value = 0.8770142695986859
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = 0.9878673343963482
$ws.Range("N17").Value = 0.7860894618181173
$ws.Range("P17").Value = 0.6661015643068852
$ws.Range("S17").Value = 0
$ws.Range("AF17").Value = 0
$ws.Range("AG17").Value = 0
$ws.Range("AT17").Value = 0.9080815959384597
$ws.Range("U18").Value = 0
$ws.Range("V18").Value = 0.7755320144143381
$ws.Range("AC18").Value = 0.6957989654682638
$ws.Range("AI18").Value = 0.9262008476325387
$ws.Range("AU18").Value = 0
$ws.Range("AX18").Value = 0
$ws.Range("Q19").Value = 0
$ws.Range("T19").Value = 0
$ws.Range("Y19").Value = 0.8036450693547653
$ws.Range("AF19").Value = 0
$ws.Range("BE19").Value = 0.9950976070420403
$ws.Range("S20").Value = 0
$ws.Range("V20").Value = 0
$ws.Range("Z20").Value = 0
$ws.Range("AA20").Value = 0.8633362626837057
$ws.Range("AF20").Value = 0
$ws.Range("AI20").Value = 0.9372112372288266
$ws.Range("AZ20").Value = 0.9907503054867941
$ws.Range("BI20").Value = 0.5215440627613782
$ws.Range("E21").Value = 0.6373157854327353
$ws.Range("R21").Value = 0
$ws.Range("W21").Value = 0.5946453920750657
$ws.Range("AF21").Value = 0.8910960817728129
$ws.Range("AX21").Value = 0.5095905866769125
$ws.Range("BD21").Value = 0.8056345597006229
$ws.Range("A22").Value = 0.6904100592492023
$ws.Range("D22").Value = 0.9158183680731049
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 0.8572906200062604
$ws.Range("J22").Value = 0
$ws.Range("R22").Value = 0.9818147104723509
$ws.Range("T22").Value = 0
$ws.Range("AD22").Value = 0.6222726418858482
$ws.Range("AX22").Value = 0
$ws.Range("BC22").Value = 0.8247767747460635
$ws.Range("BJ22").Value = 0.5773156706290681
$ws.Range("U23").Value = 0.9353900272911134
$ws.Range("X23").Value = 0
$ws.Range("AL23").Value = 0.9838537519892757
$ws.Range("AY23").Value = 0
$ws.Range("BA23").Value = 0.7528055166985219
$ws.Range("K24").Value = 0
$ws.Range("W24").Value = 0
$ws.Range("Z24").Value = 0.6007538082746021
$ws.Range("AF24").Value = 0.8346021100701619
$ws.Range("AO24").Value = 0.5600394915564415
$ws.Range("AV24").Value = 0.9126771460351077
$ws.Range("BB24").Value = 0.6559848929671289
$ws.Range("BC24").Value = 0
$ws.Range("BP24").Value = 0
$ws.Range("J25").Value = 0.907037906305221
$ws.Range("S25").Value = 0.7693458565239135
$ws.Range("AN25").Value = 0
$ws.Range("AO25").Value = 0.8351260884695965
$ws.Range("AZ25").Value = 0
$ws.Range("BF25").Value = 0
$ws.Range("T26").Value = 0
$ws.Range("X26").Value = 0.6030986292742286
$ws.Range("AB26").Value = 0
$ws.Range("AM26").Value = 0.7824565050039566
$ws.Range("AQ26").Value = 0
$ws.Range("AZ26").Value = 0
$ws.Range("BK26").Value = 0.9293965140137996
$ws.Range("C27").Value = 0
$ws.Range("J27").Value = 0.8802558731442408
$ws.Range("N27").Value = 0.8409292106810289
$ws.Range("T27").Value = 0.8684309934611976
$ws.Range("AB27").Value = 0
$ws.Range("AE27").Value = 0.7793439617804764
$ws.Range("AL27").Value = 0.6726392770039566
$ws.Range("BF27").Value = 0.7126227818120418
$ws.Range("BL27").Value = 0.8370145740539154
$ws.Range("D28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("Z28").Value = 0
$ws.Range("AA28").Value = 0
$ws.Range("AC28").Value = 0.9418149095818639
$ws.Range("AH28").Value = 0.5543875972177973
$ws.Range("BE28").Value = 0.9653548853643551
$ws.Range("BM28").Value = 0
$ws.Range("BP28").Value = 0.5862212307997183
$ws.Range("P29").Value = 0.7295576246123758
$ws.Range("R29").Value = 0.9154585480550983
$ws.Range("AB29").Value = 0.552449285939702
$ws.Range("AE29").Value = 0
$ws.Range("AP29").Value = 0.9713408183086715
$ws.Range("AS29").Value = 0.9342231653851714
$ws.Range("AV29").Value = 0
$ws.Range("BJ29").Value = 0.7731810159646919
$ws.Range("BK29").Value = 0.746066893961634
$ws.Range("BN29").Value = 0
$ws.Range("N30").Value = 0.6722111625985627
$ws.Range("O30").Value = 0
$ws.Range("P30").Value = 0.7171214935873379
$ws.Range("V30").Value = 0.7991795987167447
$ws.Range("AK30").Value = 0.797468767126637
$ws.Range("AS30").Value = 0
$ws.Range("BK30").Value = 0.985999204619721
$ws.Range("AA31").Value = 0.9907408476304629
$ws.Range("AC31").Value = 0
$ws.Range("AF31").Value = 0
$ws.Range("AG31").Value = 0.7109252558330352
$ws.Range("AP31").Value = 0
$ws.Range("AT31").Value = 0
$ws.Range("Q32").Value = 0
$ws.Range("S32").Value = 0
$ws.Range("T32").Value = 0
$ws.Range("U32").Value = 0.6549428735774866
$ws.Range("X32").Value = 0.6285374783891419
$ws.Range("AE32").Value = 0
$ws.Range("AI32").Value = 0
$ws.Range("AK32").Value = 0.719737050161015
$ws.Range("AM32").Value = 0
$ws.Range("AP32").Value = 0.9395512939609103
$ws.Range("BN32").Value = 0.5333712220120066
$ws.Range("BO32").Value = 0
$ws.Range("BP32").Value = 0.9541466999934172
$ws.Range("Q33").Value = 0
$ws.Range("AE33").Value = 0.6120843187971523
$ws.Range("AH33").Value = 0.558956708932167
$ws.Range("AI33").Value = 0.7051880782500477
$ws.Range("BK33").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("AB34").Value = 0.5802046358409273
$ws.Range("AG34").Value = 0.5961657807875291
$ws.Range("AQ34").Value = 0.7219267777656984
$ws.Range("AU34").Value = 0
$ws.Range("BD34").Value = 0.6396627761422324
$ws.Range("BL34").Value = 0.8768123540167996
$ws.Range("BM34").Value = 0.8914877918983211
$ws.Range("G35").Value = 0
$ws.Range("R35").Value = 0.8506921244746841
$ws.Range("T35").Value = 0.7053641908365541
$ws.Range("AF35").Value = 0
$ws.Range("AG35").Value = 0.6149532210964404
$ws.Range("AK35").Value = 0
$ws.Range("AM35").Value = 0.7970086882241979
$ws.Range("AU35").Value = 0.6492984835957807
$ws.Range("M36").Value = 0.7974748034062982
$ws.Range("AL36").Value = 0
$ws.Range("AS36").Value = 0.895264851037083
$ws.Range("AW36").Value = 0.6831236950475753
$ws.Range("BA36").Value = 0
$ws.Range("BF36").Value = 0
$ws.Range("BH36").Value = 0
$ws.Range("K37").Value = 0.7129824691462602
$ws.Range("M37").Value = 0.7224832187076162
$ws.Range("AD37").Value = 0.978955658870079
$ws.Range("AF37").Value = 0.8302706831266353
$ws.Range("AI37").Value = 0
$ws.Range("AN37").Value = 0
$ws.Range("AO37").Value = 0.623848621874427
$ws.Range("AP37").Value = 0
$ws.Range("BH37").Value = 0
$ws.Range("L38").Value = 0.7477308279009871
$ws.Range("M38").Value = 0
$ws.Range("W38").Value = 0.9133105953176311
$ws.Range("AA38").Value = 0.6803325623267908
$ws.Range("AJ38").Value = 0
$ws.Range("AO38").Value = 0.6001987828846814
$ws.Range("BC38").Value = 0.6849327686448367
$ws.Range("A39").Value = 0.5596716150194749
$ws.Range("I39").Value = 0
$ws.Range("Z39").Value = 0.8511718392985799
$ws.Range("AF39").Value = 0
$ws.Range("AI39").Value = 0.5295177519347691
$ws.Range("AO39").Value = 0.8038804083399824
$ws.Range("AW39").Value = 0
$ws.Range("BI39").Value = 0.8422840029032507
$ws.Range("F40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("N40").Value = 0.5374264828747672
$ws.Range("Y40").Value = 0
$ws.Range("AK40").Value = 0
$ws.Range("AZ40").Value = 0.6574810332830261
$ws.Range("BA40").Value = 0.6989220804545715
$ws.Range("BD40").Value = 0
$ws.Range("H41").Value = 0
$ws.Range("X41").Value = 0.539997759968901
$ws.Range("Y41").Value = 0.5957736750140511
$ws.Range("AK41").Value = 0.5443774821909951
$ws.Range("AL41").Value = 0.832274823523184
$ws.Range("AM41").Value = 0.7833841077789805
$ws.Range("BH41").Value = 0.8626856664703858
$ws.Range("BM41").Value = 0
$ws.Range("BP41").Value = 0
$ws.Range("O42").Value = 0.5398411525255643
$ws.Range("AC42").Value = 0.7712611049293538
$ws.Range("AE42").Value = 0
$ws.Range("AF42").Value = 0.9938731005838403
$ws.Range("AK42").Value = 0
$ws.Range("BB42").Value = 0
$ws.Range("BM42").Value = 0.6072699111609084
$ws.Range("Z43").Value = 0
$ws.Range("AH43").Value = 0.6324032961264681
$ws.Range("AT43").Value = 0.5125065280816021
$ws.Range("AY43").Value = 0.524856931081474
$ws.Range("BJ43").Value = 0
$ws.Range("BK43").Value = 0
$ws.Range("H44").Value = 0.8011248450901933
$ws.Range("AS44").Value = 0.8288563107494792
$ws.Range("AW44").Value = 0
$ws.Range("BD44").Value = 0.9188607822808941
$ws.Range("A45").Value = 0
$ws.Range("AC45").Value = 0.984923205028865
$ws.Range("AD45").Value = 0
$ws.Range("AJ45").Value = 0.7248571461293754
$ws.Range("AR45").Value = 0.9973324291375358
$ws.Range("AT45").Value = 0
$ws.Range("AX45").Value = 0
$ws.Range("BA45").Value = 0
$ws.Range("BC45").Value = 0.771376852364839
$ws.Range("K46").Value = 0
$ws.Range("M46").Value = 0.9559902397727256
$ws.Range("Q46").Value = 0.5632634736941918
$ws.Range("AE46").Value = 0
$ws.Range("AQ46").Value = 0.8052203144487925
$ws.Range("AS46").Value = 0
$ws.Range("BD46").Value = 0.8408622530773276
$ws.Range("J47").Value = 0.7273054534316707
$ws.Range("M47").Value = 0
$ws.Range("R47").Value = 0
$ws.Range("AH47").Value = 0
$ws.Range("AI47").Value = 0.7525096535358542
$ws.Range("AW47").Value = 0
$ws.Range("BB47").Value = 0.6898368446336318
$ws.Range("BO47").Value = 0
$ws.Range("I48").Value = 0.5710624044742808
$ws.Range("L48").Value = 0.5222398541355155
$ws.Range("X48").Value = 0.9267583347333141
$ws.Range("AC48").Value = 0
$ws.Range("BB48").Value = 0.6300011398120724
$ws.Range("BK48").Value = 0.9934374944143383
$ws.Range("D49").Value = 0.5312805219127856
$ws.Range("O49").Value = 0.5134763726936035
$ws.Range("AJ49").Value = 0.6369249564528694
$ws.Range("AM49").Value = 0
$ws.Range("AR49").Value = 0
$ws.Range("AU49").Value = 0
$ws.Range("AY49").Value = 0
$ws.Range("B50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("R50").Value = 0
$ws.Range("U50").Value = 0.6554390417227881
$ws.Range("V50").Value = 0
$ws.Range("AS50").Value = 0
$ws.Range("AZ50").Value = 0
$ws.Range("BF50").Value = 0.5865217695508826
$ws.Range("G51").Value = 0.9236842840667859
$ws.Range("L51").Value = 0
$ws.Range("W51").Value = 0
$ws.Range("AQ51").Value = 0.9949460703071428
$ws.Range("AW51").Value = 0
$ws.Range("BD51").Value = 0
$ws.Range("BG51").Value = 0
$ws.Range("BN51").Value = 0.9182786637699056
$ws.Range("C52").Value = 0.5764958695918487
$ws.Range("K52").Value = 0
$ws.Range("M52").Value = 0
$ws.Range("T52").Value = 0.5793461987122679
$ws.Range("Y52").Value = 0
$ws.Range("Z52").Value = 0
$ws.Range("AN52").Value = 0.537911625381654
$ws.Range("AX52").Value = 0
$ws.Range("BF52").Value = 0.5292446678610645
$ws.Range("F53").Value = 0
$ws.Range("W53").Value = 0.5672116253879316
$ws.Range("AJ53").Value = 0
$ws.Range("AN53").Value = 0.5613840201076347
$ws.Range("AS53").Value = 0
$ws.Range("BJ53").Value = 0
$ws.Range("BP53").Value = 0.6901712038279453
$ws.Range("A54").Value = 0.5427576027819816
$ws.Range("B54").Value = 0
$ws.Range("C54").Value = 0.9980634168735903
$ws.Range("M54").Value = 0
$ws.Range("X54").Value = 0.616904104216184
$ws.Range("AP54").Value = 0
$ws.Range("AU54").Value = 0.8163446256834801
$ws.Range("AV54").Value = 0.8277619189737022
$ws.Range("BM54").Value = 0
$ws.Range("BO54").Value = 0
$ws.Range("V55").Value = 0.8158480445107281
$ws.Range("X55").Value = 0
$ws.Range("AL55").Value = 0.6333331349250695
$ws.Range("AS55").Value = 0.5942571329843307
$ws.Range("BE55").Value = 0
$ws.Range("BH55").Value = 0.6674699296097912
$ws.Range("L56").Value = 0
$ws.Range("U56").Value = 0.8095952143364695
$ws.Range("AH56").Value = 0.5201591001583445
$ws.Range("AN56").Value = 0
$ws.Range("AR56").Value = 0.9079470771325111
$ws.Range("AT56").Value = 0.513629436676297
$ws.Range("AY56").Value = 0
$ws.Range("BF56").Value = 0
$ws.Range("S57").Value = 0.5441966355818356
$ws.Range("AB57").Value = 0.5681238231271883
$ws.Range("BC57").Value = 0
$ws.Range("BF57").Value = 0.8543333095355619
$ws.Range("BG57").Value = 0
$ws.Range("F58").Value = 0.6066447027529756
$ws.Range("N58").Value = 0.9514836448704627
$ws.Range("O58").Value = 0
$ws.Range("Y58").Value = 0
$ws.Range("AA58").Value = 0.7651564772686474
$ws.Range("AJ58").Value = 0
$ws.Range("AX58").Value = 0.7110338226825264
$ws.Range("AZ58").Value = 0.811658484332723
$ws.Range("BD58").Value = 0
$ws.Range("BE58").Value = 0.7122659161524869
$ws.Range("B59").Value = 0
$ws.Range("H59").Value = 0.5075719384979074
$ws.Range("P59").Value = 0
$ws.Range("AY59").Value = 0
$ws.Range("BE59").Value = 0
$ws.Range("BH59").Value = 0.7585961402234003
$ws.Range("BP59").Value = 0
$ws.Range("F60").Value = 0.8182926898483598
$ws.Range("AJ60").Value = 0
$ws.Range("AK60").Value = 0
$ws.Range("AO60").Value = 0.6361821929562715
$ws.Range("BC60").Value = 0.7072201137364006
$ws.Range("BG60").Value = 0.9844915729503094
$ws.Range("BI60").Value = 0
$ws.Range("BJ60").Value = 0
$ws.Range("BO60").Value = 0.6761380944802304
$ws.Range("J61").Value = 0
$ws.Range("T61").Value = 0.804140022614874
$ws.Range("AM61").Value = 0.6372248890928518
$ws.Range("BH61").Value = 0
$ws.Range("BN61").Value = 0
$ws.Range("A62").Value = 0.9811832898122715
$ws.Range("O62").Value = 0
$ws.Range("V62").Value = 0.8169707165598475
$ws.Range("AC62").Value = 0.8924515068029317
$ws.Range("AQ62").Value = 0
$ws.Range("BA62").Value = 0
$ws.Range("BH62").Value = 0
$ws.Range("C63").Value = 0
$ws.Range("E63").Value = 0.859657225118114
$ws.Range("Z63").Value = 0.69122588645917
$ws.Range("AC63").Value = 0.6744311235914604
$ws.Range("AD63").Value = 0.5991223193734749
$ws.Range("AG63").Value = 0
$ws.Range("AQ63").Value = 0
$ws.Range("AV63").Value = 0.6570578979665613
$ws.Range("BL63").Value = 0
$ws.Range("BO63").Value = 0.8640100302987865
$ws.Range("E64").Value = 0
$ws.Range("AA64").Value = 0.6246167884192633
$ws.Range("AH64").Value = 0.9502391270824189
$ws.Range("BK64").Value = 0
$ws.Range("BN64").Value = 0.5006198868295505
$ws.Range("BO64").Value = 0
$ws.Range("AB65").Value = 0
$ws.Range("AH65").Value = 0.858257831480457
$ws.Range("AO65").Value = 0
$ws.Range("AP65").Value = 0.736434437320975
$ws.Range("BB65").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("O66").Value = 0
$ws.Range("P66").Value = 0
$ws.Range("AC66").Value = 0
$ws.Range("AF66").Value = 0.9146916832803995
$ws.Range("AY66").Value = 0.8322296516813797
$ws.Range("BI66").Value = 0
$ws.Range("BL66").Value = 0.8729565989974424
$ws.Range("BO66").Value = 0.7545488908651033
$ws.Range("B67").Value = 0.6349191490330672
$ws.Range("F67").Value = 0.7141017167182757
$ws.Range("J67").Value = 0
$ws.Range("O67").Value = 0
$ws.Range("AF67").Value = 0
$ws.Range("AU67").Value = 0
$ws.Range("BB67").Value = 0
$ws.Range("BH67").Value = 0.9278321156545404
$ws.Range("BK67").Value = 0.7919837413526607
$ws.Range("BL67").Value = 0
$ws.Range("BN67").Value = 0.7259795554440979
$ws.Range("B68").Value = 0
$ws.Range("M68").Value = 0.6451493413999321
$ws.Range("P68").Value = 0.6299758659645782
$ws.Range("X68").Value = 0
$ws.Range("AB68").Value = 0.5379602505547524
$ws.Range("AF68").Value = 0.5196670610384806
$ws.Range("AO68").Value = 0
$ws.Range("BA68").Value = 0.9479975239778473
$ws.Range("BG68").Value = 0
